$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a "last changed" date serial that gets
# bumped by one day on every automatic refresh: 46060 -> 46061
# (2026-02-07 -> 2026-02-08), for every data row (C2:C227).
$ws.Range("C2:C227").Value = 46061
